$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (Biomarker columns shift right,
# from C:G to D:H) to make room for the new "Sex" variable.
$ws.Columns("C:C").Insert() | Out-Null

# Header for the newly inserted column.
$ws.Range("C1").Value = "Sex"

# Fill the new column with alternating Male/Female values for the data
# rows (2-26), matching the existing Patient/Control row pattern.
for ($r = 2; $r -le 26; $r++) {
    if ($r % 2 -eq 0) {
        $ws.Cells.Item($r, 3).Value = "Male"
    } else {
        $ws.Cells.Item($r, 3).Value = "Female"
    }
}

# Leave the selection where the author left it after the edit.
[void]$ws.Range("D1").Select()
